{"js": "// Change the document title from \"Version 2.\" to \"Version 1.\"\n// Original runs: \"Versi\" | \"on\" | (spellEnd) | \" 2\" | (bookmark) | \".\"\n// Target runs:   \"Version\" | (spellEnd) | \" 1.\" | (bookmark)\n\nconst body = context.document.body;\n\n// Step 1: merge the split \"Versi\"/\"on\" runs into a single \"Version\" run by\n// replacing the (cross-run) match with identical text.\nconst versionMatches = body.search(\"Version\", { matchCase: true });\nversionMatches.load(\"items\");\nawait context.sync();\nif (versionMatches.items.length > 0) {\n  versionMatches.items[0].insertText(\"Version\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Step 2: bump the version number, turning \" 2\" into \" 1.\" (the period gets\n// folded in here so the trailing standalone \".\" run becomes redundant).\nconst numberMatches = body.search(\" 2\", { matchCase: true });\nnumberMatches.load(\"items\");\nawait context.sync();\nif (numberMatches.items.length > 0) {\n  numberMatches.items[0].insertText(\" 1.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Step 3: the text now reads \"Version 1..\" (old trailing \".\" run is still\n// there) - collapse the doubled period down to a single one, which removes\n// that now-empty trailing run.\nconst dupDotMatches = body.search(\"..\", { matchCase: true });\ndupDotMatches.load(\"items\");\nawait context.sync();\nif (dupDotMatches.items.length > 0) {\n  dupDotMatches.items[0].insertText(\".\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Change the document title from \"Version 2.\" to \"Version 1.\"\n# Original runs: \"Versi\" | \"on\" | (spellEnd) | \" 2\" | (bookmark _GoBack) | \".\"\n# Target runs:   \"Version\" | (spellEnd) | \" 1.\" | (bookmark _GoBack)\n\n$wdReplaceAll   = 2\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n# Step 1: merge the split \"Versi\"/\"on\" runs into a single \"Version\" run by\n# replacing the (cross-run) match with identical text.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"Version\"\n$find1.Replacement.Text = \"Version\"\n$find1.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$wdFindContinue, [ref]$null, [ref]$null, [ref]$null, [ref]$wdReplaceAll)\n\n# Step 2: bump the version number, turning \" 2\" into \" 1.\" (the period gets\n# folded in here so the old trailing standalone \".\" run becomes redundant).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \" 2\"\n$find2.Replacement.Text = \" 1.\"\n$find2.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$wdFindContinue, [ref]$null, [ref]$null, [ref]$null, [ref]$wdReplaceAll)\n\n# Step 3: remove the now-redundant trailing \".\" run. We locate it via the\n# \"_GoBack\" bookmark (which sits right before it) instead of a text search,\n# because a Find/Replace that spans the bookmark would delete the bookmark\n# along with the matched text.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $bm = $d.Bookmarks(\"_GoBack\")\n    $afterBookmark = $bm.Range.End\n    $para = $d.Paragraphs.Item(1)\n    $paraContentEnd = $para.Range.End - 1   # exclude the paragraph mark itself\n\n    if ($paraContentEnd -gt $afterBookmark) {\n        $trailing = $d.Range($afterBookmark, $paraContentEnd)\n        $trailing.Delete()\n    }\n}\n"}
